$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The blank separator row directly beneath the title/subtitle block (row 3,
# merged into A2:B3) is being removed. Deleting it shifts every row below it
# up by one, which is exactly the row re-numbering seen in the target sheet
# (old row 4 "PROPERTY INFORMATION" -> new row 3, ..., old row 19
# "Project Status" -> new row 18). All per-row styles/merges/content move
# with their rows automatically.
$ws.Rows(3).Delete()

# Restore the cursor/selection to where the author last left it after the
# edit.
[void]$ws.Range("E10").Select()
